$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2431.7576
$ws.Range("I86").Value = 1932.1875
$ws.Range("J86").Value = 2901.9412
$ws.Range("K86").Value = 1932.1875
$ws.Range("L86").Value = 2901.9412
$ws.Range("M86").Value = -809.1875
$ws.Range("N86").Value = -5147.9412

$ws.Range("H89").Value = 2431.7576
$ws.Range("I89").Value = 1932.1875
$ws.Range("J89").Value = 2901.9412
$ws.Range("K89").Value = 9660.9375
$ws.Range("L89").Value = 14509.706
$ws.Range("M89").Value = -4044.9375
$ws.Range("N89").Value = -25741.706

$ws.Range("H106").Value = 3836.6924
$ws.Range("I106").Value = 3730.7778
$ws.Range("J106").Value = 4075
$ws.Range("K106").Value = 3730.7778
$ws.Range("L106").Value = 4075
$ws.Range("M106").Value = -3099.7778
$ws.Range("N106").Value = -5337

$ws.Range("H137").Value = 888.37933
$ws.Range("I137").Value = 841.25
$ws.Range("J137").Value = 946.38464
$ws.Range("K137").Value = 2523.75
$ws.Range("L137").Value = 2839.15392
$ws.Range("M137").Value = 26.25
$ws.Range("N137").Value = -7939.15392

$ws.Range("H138").Value = 2070.56
$ws.Range("I138").Value = 1006.25
$ws.Range("J138").Value = 2571.4119
$ws.Range("K138").Value = 3018.75
$ws.Range("L138").Value = 7714.2357
$ws.Range("M138").Value = 2121.25
$ws.Range("N138").Value = -17994.2357

$ws.Range("H141").Value = 2691.5833
$ws.Range("I141").Value = 1110.4286
$ws.Range("J141").Value = 8225.625
$ws.Range("K141").Value = 3331.2858
$ws.Range("L141").Value = 24676.875
$ws.Range("M141").Value = 1848.7142
$ws.Range("N141").Value = -35036.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20001.01
$ws.Range("I32").Value = 19764.78
$ws.Range("J32").Value = 20867.191
$ws.Range("K32").Value = 19764.78
$ws.Range("L32").Value = 20867.191
$ws.Range("M32").Value = -19477.78
$ws.Range("N32").Value = -21441.191

$ws.Range("H74").Value = 1061.25
$ws.Range("I74").Value = 1030.0571
$ws.Range("J74").Value = 1279.6
$ws.Range("K74").Value = 1030.0571
$ws.Range("L74").Value = 1279.6
$ws.Range("M74").Value = -156.0571
$ws.Range("N74").Value = -3027.6

$ws.Range("H77").Value = 1061.25
$ws.Range("I77").Value = 1030.0571
$ws.Range("J77").Value = 1279.6
$ws.Range("K77").Value = 5150.2855
$ws.Range("L77").Value = 6398
$ws.Range("M77").Value = -782.2855
$ws.Range("N77").Value = -15134

$ws.Range("H132").Value = 2177.68
$ws.Range("I132").Value = 2062.739
$ws.Range("J132").Value = 3499.5
$ws.Range("K132").Value = 6188.217000000001
$ws.Range("L132").Value = 10498.5
$ws.Range("M132").Value = -3658.217000000001
$ws.Range("N132").Value = -15558.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 14924.373
$ws.Range("I134").Value = 1378.2742
$ws.Range("J134").Value = 79528.84
$ws.Range("K134").Value = 4134.8226
$ws.Range("L134").Value = 238586.52
$ws.Range("M134").Value = -1599.8226
$ws.Range("N134").Value = -243656.52

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1917291.6
$ws.Range("I31").Value = 3002631.2
$ws.Range("J31").Value = 1986.1765
$ws.Range("K31").Value = 3002631.2
$ws.Range("L31").Value = 1986.1765
$ws.Range("M31").Value = -3002336.2
$ws.Range("N31").Value = -2576.1765

$ws.Range("H34").Value = 1917291.6
$ws.Range("I34").Value = 3002631.2
$ws.Range("J34").Value = 1986.1765
$ws.Range("K34").Value = 3002631.2
$ws.Range("L34").Value = 1986.1765
$ws.Range("M34").Value = -3002429.2
$ws.Range("N34").Value = -2390.1765

$ws.Range("H134").Value = 1214.619
$ws.Range("I134").Value = 1173.5682
$ws.Range("J134").Value = 1309.6842
$ws.Range("K134").Value = 3520.7046
$ws.Range("L134").Value = 3929.0526
$ws.Range("M134").Value = -985.7046
$ws.Range("N134").Value = -8999.052599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 1091.6
$ws.Range("I45").Value = 716.6667
$ws.Range("J45").Value = 1252.2858
$ws.Range("K45").Value = 2150.0001
$ws.Range("L45").Value = 3756.8574
$ws.Range("M45").Value = -1618.0001
$ws.Range("N45").Value = -4820.857400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1136.8462
$ws.Range("I102").Value = 1256.3
$ws.Range("J102").Value = 738.6667
$ws.Range("K102").Value = 1256.3
$ws.Range("L102").Value = 738.6667
$ws.Range("M102").Value = 365.7
$ws.Range("N102").Value = -3982.6667

$ws.Range("H132").Value = 2307.5366
$ws.Range("I132").Value = 2202.8696
$ws.Range("J132").Value = 2441.2778
$ws.Range("K132").Value = 6608.6088
$ws.Range("L132").Value = 7323.8334
$ws.Range("M132").Value = -4078.6088
$ws.Range("N132").Value = -12383.8334

$ws.Range("H136").Value = 9845
$ws.Range("J136").Value = 9845
$ws.Range("L136").Value = 29535
$ws.Range("N136").Value = -34635

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5572.222
$ws.Range("I16").Value = 8028.625
$ws.Range("K16").Value = 8028.625
$ws.Range("M16").Value = -7858.625

$ws.Range("H127").Value = 54781.25
$ws.Range("J127").Value = 54781.25
$ws.Range("L127").Value = 54781.25
$ws.Range("N127").Value = -64701.25

$ws.Range("H136").Value = 1851.6666
$ws.Range("I136").Value = 1039.4048
$ws.Range("J136").Value = 5642.222
$ws.Range("K136").Value = 3118.2144
$ws.Range("L136").Value = 16926.666
$ws.Range("M136").Value = -568.2143999999998
$ws.Range("N136").Value = -22026.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1656.5714
$ws.Range("I81").Value = 1480.5333
$ws.Range("J81").Value = 2096.6667
$ws.Range("K81").Value = 2961.0666
$ws.Range("L81").Value = 4193.3334
$ws.Range("M81").Value = -1900.0666
$ws.Range("N81").Value = -6315.3334

$ws.Range("H84").Value = 1656.5714
$ws.Range("I84").Value = 1480.5333
$ws.Range("J84").Value = 2096.6667
$ws.Range("K84").Value = 14805.333
$ws.Range("L84").Value = 20966.667
$ws.Range("M84").Value = -9501.333000000001
$ws.Range("N84").Value = -31574.667

$ws.Range("H107").Value = 760
$ws.Range("I107").Value = 951.13336
$ws.Range("J107").Value = 401.625
$ws.Range("K107").Value = 2853.40008
$ws.Range("L107").Value = 1204.875
$ws.Range("M107").Value = -933.4000800000003
$ws.Range("N107").Value = -5044.875

